$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing shared string text (B9 cell content)
$ws.Range("B9").Value = "Redoing backend for ui controllers so its easier to use. And some frontend for the movement system."

# Add hours value for row 9
$ws.Range("C9").Value = 6

# Add new row 10 with a date in A10, formatted same as other date cells (copy format from A9)
$ws.Range("A9").Copy() | Out-Null
$ws.Range("A10").PasteSpecial(-4122) | Out-Null  # xlPasteFormats
$ws.Range("A10").Value = 45624

# Update selection to match diff (B13)
$ws.Range("B13").Select() | Out-Null
